$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Functions")

# Set the "Written" column (B) to "X" for rows 120, 121, 331, 373
$ws.Cells.Item(120, 2).Value = "X"
$ws.Cells.Item(121, 2).Value = "X"
$ws.Cells.Item(331, 2).Value = "X"
$ws.Cells.Item(373, 2).Value = "X"

# Hide rows 89, 120, 121, 331, 373
$ws.Rows.Item(89).Hidden = $true
$ws.Rows.Item(120).Hidden = $true
$ws.Rows.Item(121).Hidden = $true
$ws.Rows.Item(331).Hidden = $true
$ws.Rows.Item(373).Hidden = $true

# Update the active selection to A269
$ws.Range("A269").Select()
